# "Add files via upload" — append a new date column (23-ago) to the
# "Dataframe Fam" weekly-family table: a new header cell in row 1 plus one
# count per product row (rows 2-11), mirroring the existing AZ ("22-ago")
# column's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AZ = 52 (existing last data column, "22-ago").
# Column BA = 53 (new column being appended, "23-ago").
$newCol = 53

# New counts for the "23-ago" column, one per data row (rows 2-11).
$values = @{
    2  = 15
    3  = 13
    4  = 10
    5  = 12
    6  = 12
    7  = 12
    8  = 13
    9  = 17
    10 = 14
    11 = 15
}

# Match the formatting of the existing data column (AZ) before writing
# values, so the new column's style indices line up with the rest of the
# table instead of minting new ones.
$dataRange = $ws.Range($ws.Cells.Item(2, $newCol), $ws.Cells.Item(11, $newCol))
$dataRange.HorizontalAlignment = -4108
$dataRange.NumberFormat = "0"
$ws.Cells.Item(1, $newCol).NumberFormat = "@"

# Header cell: new date label for the appended column.
$ws.Cells.Item(1, $newCol).Value = "23-ago"

# Data cells: one count per product row.
foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, $newCol).Value = $values[$row]
}

# Move/restore the active selection the way it landed after the edit.
$ws.Range("BA12").Select()
